$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7916129955631771
$ws.Range("C2").Value = -1.579792844692939
$ws.Range("D2").Value = -0.9155274120541215
$ws.Range("E2").Value = -0.1460550576300293
$ws.Range("F2").Value = -0.2386646425575917
$ws.Range("G2").Value = 0.4822622889884105
$ws.Range("H2").Value = -0.2266896356254442
$ws.Range("I2").Value = 0.2654163566970329
$ws.Range("J2").Value = 0.7352211039700201
$ws.Range("K2").Value = 1.048046298935466
$ws.Range("B3").Value = -3.727363316492332
$ws.Range("C3").Value = -3.063097883853514
$ws.Range("D3").Value = -2.293625529429422
$ws.Range("E3").Value = -2.386235114356984
$ws.Range("F3").Value = -1.665308182810982
$ws.Range("G3").Value = -2.374260107424837
$ws.Range("H3").Value = -1.88215411510236
$ws.Range("I3").Value = -1.412349367829372
$ws.Range("J3").Value = -1.099524172863926
$ws.Range("K3").Value = -1.977926754115506
$ws.Range("B4").Value = 0.376932102669816
$ws.Range("C4").Value = 1.146404457093908
$ws.Range("D4").Value = 1.053794872166346
$ws.Range("E4").Value = 1.774721803712348
$ws.Range("F4").Value = 1.065769879098493
$ws.Range("G4").Value = 1.55787587142097
$ws.Range("H4").Value = 2.027680618693958
$ws.Range("I4").Value = 2.340505813659403
$ws.Range("J4").Value = 1.462103232407824
$ws.Range("K4").Value = 0.7321912172293545
$ws.Range("B5").Value = 1.207578635508109
$ws.Range("C5").Value = 1.114969050580547
$ws.Range("D5").Value = 1.835895982126549
$ws.Range("E5").Value = 1.126944057512694
$ws.Range("F5").Value = 1.619050049835171
$ws.Range("G5").Value = 2.088854797108159
$ws.Range("H5").Value = 2.401679992073604
$ws.Range("I5").Value = 1.523277410822025
$ws.Range("J5").Value = 0.7933653956435556
$ws.Range("K5").Value = 2.455544492033183
$ws.Range("B6").Value = -0.9264868865757077
$ws.Range("C6").Value = -0.2055599550297054
$ws.Range("D6").Value = -0.9145118796435601
$ws.Range("E6").Value = -0.422405887321083
$ws.Range("F6").Value = 0.04739885995190407
$ws.Range("G6").Value = 0.36022405491735
$ws.Range("H6").Value = -0.5181785263342299
$ws.Range("I6").Value = -1.248090541512699
$ws.Range("J6").Value = 0.4140885548769285
$ws.Range("K6").Value = -0.2051307335183153
$ws.Range("B7").Value = 0.3770345820039356
$ws.Range("C7").Value = -0.3319173426099191
$ws.Range("D7").Value = 0.160188649712558
$ws.Range("E7").Value = 0.6299933969855451
$ws.Range("F7").Value = 0.942818591950991
$ws.Range("G7").Value = 0.06441601069941108
$ws.Range("H7").Value = -0.6654960044790579
$ws.Range("I7").Value = 0.9966830919105695
$ws.Range("J7").Value = 0.3774638035153257
$ws.Range("K7").Value = 0.6147675671350392
$ws.Range("B8").Value = -0.4275923834192769
$ws.Range("C8").Value = 0.0645136089032002
$ws.Range("D8").Value = 0.5343183561761873
$ws.Range("E8").Value = 0.8471435511416332
$ws.Range("F8").Value = -0.03125903010994671
$ws.Range("G8").Value = -0.7611710452884157
$ws.Range("H8").Value = 0.9010080511012117
$ws.Range("I8").Value = 0.2817887627059679
$ws.Range("J8").Value = 0.5190925263256815
$ws.Range("K8").Value = 0.4158151564502698
$ws.Range("B9").Value = 0.324932645901923
$ws.Range("C9").Value = 0.7947373931749101
$ws.Range("D9").Value = 1.107562588140356
$ws.Range("E9").Value = 0.2291600068887761
$ws.Range("F9").Value = -0.5007520082896928
$ws.Range("G9").Value = 1.161427088099934
$ws.Range("H9").Value = 0.5422077997046907
$ws.Range("I9").Value = 0.7795115633244043
$ws.Range("J9").Value = 0.6762341934489926
$ws.Range("K9").Value = -0.1153642338804421
$ws.Range("B10").Value = -0.04071760298358112
$ws.Range("C10").Value = 0.2721075919818648
$ws.Range("D10").Value = -0.6062949892697151
$ws.Range("E10").Value = -1.336207004448184
$ws.Range("F10").Value = 0.3259720919414433
$ws.Range("G10").Value = -0.2932471964538005
$ws.Range("H10").Value = -0.05594343283408693
$ws.Range("I10").Value = -0.1592208027094986
$ws.Range("J10").Value = -0.9508192300389333
$ws.Range("K10").Value = -0.3403303223714723
$ws.Range("B11").Value = 0.3721869518844864
$ws.Range("C11").Value = -0.5062156293670936
$ws.Range("D11").Value = -1.236127644545562
$ws.Range("E11").Value = 0.4260514518440648
$ws.Range("F11").Value = -0.193167836551179
$ws.Range("G11").Value = 0.04413592706853459
$ws.Range("H11").Value = -0.05914144280687711
$ws.Range("I11").Value = -0.8507398701363118
$ws.Range("J11").Value = -0.2402509624688508
$ws.Range("K11").Value = -0.4017729932881683
$ws.Range("B12").Value = -0.1524291232873974
$ws.Range("C12").Value = -0.8823411384658664
$ws.Range("D12").Value = 0.779837957923761
$ws.Range("E12").Value = 0.1606186695285172
$ws.Range("F12").Value = 0.3979224331482308
$ws.Range("G12").Value = 0.2946450632728191
$ws.Range("H12").Value = -0.4969533640566156
$ws.Range("I12").Value = 0.1135355436108454
$ws.Range("J12").Value = -0.04798648720847212
$ws.Range("B13").Value = -1.030518528898312
$ws.Range("C13").Value = 0.6316605674913157
$ws.Range("D13").Value = 0.0124412790960719
$ws.Range("E13").Value = 0.2497450427157855
$ws.Range("F13").Value = 0.1464676728403738
$ws.Range("G13").Value = -0.6451307544890609
$ws.Range("H13").Value = -0.03464184682159993
$ws.Range("I13").Value = -0.1961638776409175
$ws.Range("B14").Value = 0.4742145784871607
$ws.Range("C14").Value = -0.1450047099080831
$ws.Range("D14").Value = 0.0922990537116305
$ws.Range("E14").Value = -0.0109783161637812
$ws.Range("F14").Value = -0.8025767434932158
$ws.Range("G14").Value = -0.1920878358257549
$ws.Range("H14").Value = -0.3536098666450724
$ws.Range("B15").Value = 0.3556547466179877
$ws.Range("C15").Value = 0.5929585102377013
$ws.Range("D15").Value = 0.4896811403622896
$ws.Range("E15").Value = -0.3019172869671451
$ws.Range("F15").Value = 0.3085716207003159
$ws.Range("G15").Value = 0.1470495898809984
$ws.Range("B16").Value = 0.3126006297022321
$ws.Range("C16").Value = 0.2093232598268204
$ws.Range("D16").Value = -0.5822751675026142
$ws.Range("E16").Value = 0.02821374016484672
$ws.Range("F16").Value = -0.1333082906544708
$ws.Range("B17").Value = 0.3812981176718321
$ws.Range("C17").Value = -0.4103003096576026
$ws.Range("D17").Value = 0.2001885980098584
$ws.Range("E17").Value = 0.03866656719054083
$ws.Range("B18").Value = -0.716162849403934
$ws.Range("C18").Value = -0.1056739417364731
$ws.Range("D18").Value = -0.2671959725557906
$ws.Range("B19").Value = 0.506656010950813
$ws.Range("C19").Value = 0.3451339801314955
$ws.Range("B20").Value = -0.343237405067616
